# The code-review checklist table's header row holds "Sprint No." and
# "Review Date" values that need updating for this sprint's review.
$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# "Sprint No." value: "1" -> "2"  (Table 1, Row 2, Col 4)
$sprintCell = $table.Cell(2, 4)
$sprintCell.Range.Find.Execute("1", $true, $false, $false, $false, $false, `
    $true, 0, $false, "2", 1) | Out-Null

# "Review Date" value: "02/09/18" -> "02/21/18"  (Table 1, Row 3, Col 2)
$dateCell = $table.Cell(3, 2)
$dateCell.Range.Find.Execute("02/09/18", $true, $false, $false, $false, $false, `
    $true, 0, $false, "02/21/18", 1) | Out-Null
